$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-case description text in column M (rows 2-27):
# version reference changes from V1.65 to V1.64
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 13)  # column M = 13
    $cell.Value = "製作依據之需求規格書與版本：PJ201800012_URS_5管理性作業_V1.64.DOCX"
}

# Update the active selection on the sheet from R2 to A2
$ws.Range("A2").Select()
